# Applies the "output generated at 456a3b4" update to 广州-漫展信息.xlsx
# Sheets (tab order): 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet1 (展览): "想去人数" (F) counter bumps, plus one "最低票价" (G40) fix ---
$ws1.Range("F7").Value  = 2642
$ws1.Range("F9").Value  = 1292
$ws1.Range("F11").Value = 624
$ws1.Range("F12").Value = 935
$ws1.Range("F14").Value = 295
$ws1.Range("F16").Value = 740
$ws1.Range("F17").Value = 788
$ws1.Range("F18").Value = 225
$ws1.Range("F19").Value = 525
$ws1.Range("F20").Value = 1137
$ws1.Range("F21").Value = 103
$ws1.Range("F22").Value = 633
$ws1.Range("F23").Value = 608
$ws1.Range("F24").Value = 225
$ws1.Range("F25").Value = 315
$ws1.Range("F26").Value = 312
$ws1.Range("F27").Value = 693
$ws1.Range("F28").Value = 569
$ws1.Range("F29").Value = 5668
$ws1.Range("F30").Value = 5668
$ws1.Range("F31").Value = 497
$ws1.Range("F32").Value = 44
$ws1.Range("F34").Value = 32
$ws1.Range("F35").Value = 174
$ws1.Range("F36").Value = 1637
$ws1.Range("F37").Value = 11
$ws1.Range("F38").Value = 91
$ws1.Range("F39").Value = 446
$ws1.Range("G40").Value = 68
$ws1.Range("F42").Value = 152
$ws1.Range("F43").Value = 11
$ws1.Range("F44").Value = 75
$ws1.Range("F46").Value = 143
$ws1.Range("F48").Value = 118

# --- Sheet2 (演出) ---
$ws2.Range("F3").Value  = 86
$ws2.Range("F9").Value  = 41
$ws2.Range("F13").Value = 4411
$ws2.Range("F14").Value = 39
$ws2.Range("F17").Value = 37
$ws2.Range("G17").Value = 180
$ws2.Range("F18").Value = 209

# --- Sheet3 (本地生活): remove the expired "天官赐福" flash-store entry (row 2);
#     the "全职高手 x HAPPY ZOO" cafe entry (old row 3) shifts up to row 2, and
#     its "想去人数" count is refreshed from 740 -> 744.
$ws3.Rows.Item(2).Delete()
$ws3.Range("A2").Value = 1
$ws3.Range("F2").Value = 744
$ws3.Range("G2").Value = 10

# --- Sheet4 (全部类型) ---
$ws4.Range("F3").Value  = 86
$ws4.Range("F5").Value  = 2642
$ws4.Range("F7").Value  = 1292
$ws4.Range("F9").Value  = 624
$ws4.Range("F10").Value = 935
$ws4.Range("F12").Value = 295
$ws4.Range("F15").Value = 740
$ws4.Range("F18").Value = 788
$ws4.Range("F19").Value = 225
$ws4.Range("F20").Value = 525
$ws4.Range("F21").Value = 1137
$ws4.Range("F24").Value = 103
$ws4.Range("F25").Value = 41
$ws4.Range("F26").Value = 633
$ws4.Range("F27").Value = 608
$ws4.Range("F28").Value = 315
$ws4.Range("F29").Value = 312
$ws4.Range("F30").Value = 569
$ws4.Range("F31").Value = 5668
$ws4.Range("F33").Value = 497
$ws4.Range("F34").Value = 44
$ws4.Range("F35").Value = 32
$ws4.Range("F36").Value = 174
$ws4.Range("F37").Value = 1637
$ws4.Range("F38").Value = 11
$ws4.Range("F39").Value = 446
$ws4.Range("F40").Value = 39
$ws4.Range("F41").Value = 39
$ws4.Range("G42").Value = 68
$ws4.Range("F44").Value = 37
$ws4.Range("G44").Value = 180
$ws4.Range("F45").Value = 152
$ws4.Range("F46").Value = 75
$ws4.Range("F49").Value = 118
